# Aggiornamento fino a 27/05: append new daily rows (256-269) to the
# "Marano" report sheet, extending the data range from A1:D255 to A1:D269.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
  @(44330, 0, 3, 56.81818181818181),
  @(44331, 1, 4, 75.75757575757575),
  @(44332, 1, 3, 56.81818181818181),
  @(44333, 0, 3, 56.81818181818181),
  @(44334, 0, 2, 37.87878787878788),
  @(44335, 0, 2, 37.87878787878788),
  @(44336, 0, 2, 37.87878787878788),
  @(44337, 0, 2, 37.87878787878788),
  @(44338, 0, 1, 18.93939393939394),
  @(44339, 0, 0, 0),
  @(44340, 0, 0, 0),
  @(44341, 0, 0, 0),
  @(44342, 0, 0, 0),
  @(44343, 0, 0, 0)
)

$lastRow = 255
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $data.Count

# Carry the existing row formatting (column A date style w/ border+bold+centered,
# matches the style already used by the rest of the column) down onto the new rows
# before writing the values, exactly like dragging the fill handle in Excel.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$firstNewRow`:D$lastNewRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$r = $firstNewRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
